$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 358.5
$ws.Range("I2").Value = 280.125
$ws.Range("K2").Value = 280.125
$ws.Range("M2").Value = -167.125
$ws.Range("H8").Value = 34.833332
$ws.Range("I8").Value = 34.833332
$ws.Range("K8").Value = 104.499996
$ws.Range("M8").Value = 34.500004
$ws.Range("H21").Value = 14000
$ws.Range("I21").Value = 8000
$ws.Range("K21").Value = 8000
$ws.Range("M21").Value = -7532
$ws.Range("H23").Value = 14000
$ws.Range("I23").Value = 8000
$ws.Range("K23").Value = 8000
$ws.Range("M23").Value = -7766
$ws.Range("H38").Value = 321.23077
$ws.Range("J38").Value = 954
$ws.Range("L38").Value = 2862
$ws.Range("N38").Value = -3606
$ws.Range("H42").Value = 608.3570999999999
$ws.Range("I42").Value = 101.7
$ws.Range("J42").Value = 1875
$ws.Range("K42").Value = 305.1
$ws.Range("L42").Value = 5625
$ws.Range("M42").Value = -75.10000000000002
$ws.Range("N42").Value = -6085
$ws.Range("H58").Value = 1904.4445
$ws.Range("I58").Value = 1085
$ws.Range("J58").Value = 2560
$ws.Range("K58").Value = 3255
$ws.Range("L58").Value = 7680
$ws.Range("M58").Value = -3105
$ws.Range("N58").Value = -7980
$ws.Range("H62").Value = 5007.5
$ws.Range("I62").Value = 2908.4375
$ws.Range("J62").Value = 21800
$ws.Range("K62").Value = 2908.4375
$ws.Range("L62").Value = 21800
$ws.Range("M62").Value = -2284.4375
$ws.Range("N62").Value = -23048
$ws.Range("H65").Value = 5007.5
$ws.Range("I65").Value = 2908.4375
$ws.Range("J65").Value = 21800
$ws.Range("K65").Value = 14542.1875
$ws.Range("L65").Value = 109000
$ws.Range("M65").Value = -11422.1875
$ws.Range("N65").Value = -115240
$ws.Range("H106").Value = 2916.3333
$ws.Range("I106").Value = 2916.3333
$ws.Range("K106").Value = 2916.3333
$ws.Range("M106").Value = -2285.3333
$ws.Range("H111").Value = 2962.5
$ws.Range("I111").Value = 2625
$ws.Range("K111").Value = 7875
$ws.Range("M111").Value = -4808
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7439.8774
$ws.Range("I61").Value = 4255.364
$ws.Range("J61").Value = 14007.9375
$ws.Range("K61").Value = 4255.364
$ws.Range("L61").Value = 14007.9375
$ws.Range("M61").Value = -4043.364
$ws.Range("N61").Value = -14431.9375
$ws.Range("H74").Value = 4308.974
$ws.Range("I74").Value = 1758.3429
$ws.Range("K74").Value = 1758.3429
$ws.Range("M74").Value = -884.3429000000001
$ws.Range("H77").Value = 4308.974
$ws.Range("I77").Value = 1758.3429
$ws.Range("K77").Value = 8791.7145
$ws.Range("M77").Value = -4423.7145
$ws.Range("H136").Value = 7439.8774
$ws.Range("I136").Value = 4255.364
$ws.Range("J136").Value = 14007.9375
$ws.Range("K136").Value = 12766.092
$ws.Range("L136").Value = 42023.8125
$ws.Range("M136").Value = -10216.092
$ws.Range("N136").Value = -47123.8125

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2013.8889
$ws.Range("I109").Value = 1023.75
$ws.Range("J109").Value = 2806
$ws.Range("K109").Value = 3071.25
$ws.Range("L109").Value = 8418
$ws.Range("M109").Value = -2031.25
$ws.Range("N109").Value = -10498
$ws.Range("H126").Value = 3356.2856
$ws.Range("I126").Value = 1100
$ws.Range("J126").Value = 3732.3333
$ws.Range("K126").Value = 3300
$ws.Range("L126").Value = 11196.9999
$ws.Range("M126").Value = 1640
$ws.Range("N126").Value = -21076.9999
$ws.Range("H132").Value = 1366
$ws.Range("I132").Value = 1783.8182
$ws.Range("J132").Value = 1174.5
$ws.Range("K132").Value = 16054.3638
$ws.Range("L132").Value = 10570.5
$ws.Range("M132").Value = -13524.3638
$ws.Range("N132").Value = -15630.5
$ws.Range("H134").Value = 4998.864
$ws.Range("I134").Value = 5104.9165
$ws.Range("J134").Value = 4871.6
$ws.Range("K134").Value = 15314.7495
$ws.Range("L134").Value = 14614.8
$ws.Range("M134").Value = -10244.7495
$ws.Range("N134").Value = -24754.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 7514548.5
$ws.Range("I18").Value = 30000004
$ws.Range("J18").Value = 19396.666
$ws.Range("K18").Value = 30000004
$ws.Range("L18").Value = 19396.666
$ws.Range("M18").Value = -29999711
$ws.Range("N18").Value = -19982.666
$ws.Range("H42").Value = 60467.5
$ws.Range("J42").Value = 60467.5
$ws.Range("L42").Value = 60467.5
$ws.Range("N42").Value = -61437.5
$ws.Range("H115").Value = 60467.5
$ws.Range("J115").Value = 60467.5
$ws.Range("L115").Value = 60467.5
$ws.Range("N115").Value = -62817.5
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 3750
$ws.Range("J20").Value = 3750
$ws.Range("L20").Value = 3750
$ws.Range("N20").Value = -4202
$ws.Range("H23").Value = 1000000
$ws.Range("I23").Value = 1000000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 1000000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -999770
$ws.Range("N23").ClearContents()
$ws.Range("H61").Value = 732011.2
$ws.Range("I61").Value = 28098.875
$ws.Range("J61").Value = 1670561
$ws.Range("K61").Value = 28098.875
$ws.Range("L61").Value = 1670561
$ws.Range("M61").Value = -27896.875
$ws.Range("N61").Value = -1670965
$ws.Range("H113").Value = 732011.2
$ws.Range("I113").Value = 28098.875
$ws.Range("J113").Value = 1670561
$ws.Range("K113").Value = 28098.875
$ws.Range("L113").Value = 1670561
$ws.Range("M113").Value = -25928.875
$ws.Range("N113").Value = -1674901

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 12502982
$ws.Range("I81").Value = 1616.8334
$ws.Range("J81").Value = 20003802
$ws.Range("K81").Value = 3233.6668
$ws.Range("L81").Value = 40007604
$ws.Range("M81").Value = -2172.6668
$ws.Range("N81").Value = -40009726
$ws.Range("H84").Value = 12502982
$ws.Range("I84").Value = 1616.8334
$ws.Range("J84").Value = 20003802
$ws.Range("K84").Value = 16168.334
$ws.Range("L84").Value = 200038020
$ws.Range("M84").Value = -10864.334
$ws.Range("N84").Value = -200048628
